{"js": "// The paragraph \"<id>p057v_2</id>\" is currently split across three runs:\n//   <id>   (Courier New, color 7f6000)\n//   p057v_2 (color 000000)\n//   </id>  (Courier New, color 7f6000)\n// The edit merges them into a single run containing the whole literal\n// text \"<id>p057v_2</id>\", taking on the formatting of the first run.\nconst searchResults = context.document.body.search(\"<id>p057v_2</id>\", {\n  matchCase: true,\n  matchWildcards: false\n});\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not find target text '<id>p057v_2</id>' in document body.\");\n}\n\n// Replacing the matched range with the same literal text collapses the\n// three runs it spans into one run, inheriting the first run's formatting.\nconst target = searchResults.items[0];\ntarget.insertText(\"<id>p057v_2</id>\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The paragraph \"<id>p057v_2</id>\" is currently split across three runs:\n#   <id>    (Courier New, color 7f6000)\n#   p057v_2 (color 000000)\n#   </id>   (Courier New, color 7f6000)\n# The edit merges them into a single run containing the whole literal\n# text \"<id>p057v_2</id>\", keeping the formatting of the FIRST run.\n\n$d = $word.ActiveDocument\n\n# Locate the target text anywhere in the document body.\n$r = $d.Content\n$find = $r.Find\n$find.ClearFormatting()\n$find.Text = \"<id>p057v_2</id>\"\n$found = $find.Execute()\n\nif ($found) {\n    $fullText = $r.Text\n    $start = $r.Start\n    $end = $r.End\n\n    # Determine how many characters at the start belong to the first run\n    # by walking forward until the character formatting (font name/color)\n    # changes - that's the boundary of the original first run.\n    $firstChar = $d.Range($start, $start + 1)\n    $refColor = $firstChar.Font.Color\n    $refName = $firstChar.Font.Name\n\n    $boundary = $end\n    for ($i = $start + 1; $i -lt $end; $i++) {\n        $c = $d.Range($i, $i + 1)\n        if ($c.Font.Color -ne $refColor -or $c.Font.Name -ne $refName) {\n            $boundary = $i\n            break\n        }\n    }\n\n    $firstRunLen = $boundary - $start\n\n    # Keep the original first run (and its formatting/xml:space) intact,\n    # delete the remaining runs' text, then append the rest of the target\n    # text onto that first run so everything collapses into one run.\n    $firstRunRange = $d.Range($start, $boundary)\n    $remainder = $d.Range($boundary, $end)\n    $remainder.Delete()\n    $firstRunRange.InsertAfter($fullText.Substring($firstRunLen))\n}\n"}
